# fix: add user setting sample APi seqs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = "2025-03-06 12:03:12"
$ws.Range("F2").Value = "http://47.97.114.24:5230/api/v1/resource/16"
$ws.Range("G2").Value = "/api/v1/resource/16"
$ws.Range("M2").Value = 0.003
$ws.Range("N2").Value = 0
$ws.Range("Q2").Value = $true

# --- Row 3 ---
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = "2025-03-06 12:03:12"
$ws.Range("F3").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G3").Value = "/api/v1/memo/21"
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $true

# --- Row 4 ---
$ws.Range("D4").Value = "2025-03-06 12:03:12"
$ws.Range("M4").Value = 0.002
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = $false

# --- Row 5 ---
$ws.Range("D5").Value = "2025-03-06 12:03:12"
$ws.Range("F5").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G5").Value = "/api/v1/memo/21"
$ws.Range("M5").Value = 0.003
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = $false

# --- Row 6 ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = "2025-03-06 12:03:12"
$ws.Range("F6").Value = "http://49.234.6.241:5230/api/v1/resource/16"
$ws.Range("G6").Value = "/api/v1/resource/16"
$ws.Range("M6").Value = 0.002
$ws.Range("N6").Value = 0
$ws.Range("Q6").Value = $true

# --- Remove row 7 (old last data row, now merged into the dataset above) ---
$ws.Rows.Item(7).Delete()
